$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Cachoeirinha, 2025-04-04 00:00:00.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Cachoeirinha, 04 de April de 2025.",
    2
)
